$d = $word.ActiveDocument

# 1) Insert three new "Normal" style paragraphs at the very start of the document:
#    "Ryan Cummings", "Alex Ackerlund", "Bruce Weston"
$r1 = $d.Paragraphs(1).Range
$r1.InsertBefore("Ryan Cummings`r")

$r2 = $d.Paragraphs(2).Range
$r2.InsertBefore("Alex Ackerlund`r")

$r3 = $d.Paragraphs(3).Range
$r3.InsertBefore("Bruce Weston`r")

# 2) Reword how the tool surmises design patterns from file names.
$old2 = "developers to name their files according to the design patterns they used, the tool surmises"
$new2 = "developers to choose descriptive file names, the tool surmises"
$d.Content.Find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# 3) Append the new "machine learning" discussion to the end of the paragraph that
#    ends with "A quick look at the source code confirms this intuition."
$old3 = "A quick look at the source code confirms this intuition."
$new3 = "A quick look at the source code confirms this intuition.  How we would search for design patterns would be to use machine learning.  Assuming our algorithm would have access to the code and not just the file names, we would train it to look for relationships between classes/interfaces that point to particular patterns.  If, for example, a class had an instance field corresponding to an interface and that instance field was used to dynamically change the behavior of subclasses at runtime, we could train the algorithm to see the Strategy Pattern. By recognizing tell-tale relationships with machine learning, we could accurately gauge which patterns were being used in any given source code. "
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2) | Out-Null

# 4) Replace "git clone <repo url> " with "git clone <URL> " (and drop the spell-check
#    markup that used to wrap "url").
$old4 = "git clone <repo url> "
$new4 = "git clone <URL> "
$d.Content.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# 5) "...page for future assignments," -> "...page for current and future assignments,"
$old5 = " page for future assignments, "
$new5 = " page for current and future assignments, "
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2) | Out-Null

Write-Output "done"
